$wb = $excel.ActiveWorkbook
$terms = $wb.Worksheets.Item("Terms")

# Update the Terms sheet data: swap Male/Female order, rename Other->NAP, add NKN/NST rows
$terms.Range("A3").Value = "PersonGenderEnum/Female"
$terms.Range("A4").Value = "PersonGenderEnum/Male"
$terms.Range("A5").Value = "PersonGenderEnum/NAP"
$terms.Range("A6").Value = "PersonGenderEnum/NKN"
$terms.Range("C6").Value = "xs:token"
$terms.Range("A7").Value = "PersonGenderEnum/NST"
$terms.Range("C7").Value = "xs:token"

# Add hyperlink over range B2:B7 -> Human Sex (display text is the bare URL)
$url = "https://op.europa.eu/en/web/eu-vocabularies/at-concept-scheme/-/resource/authority/human-sex/?target=Browse&uri=http://publications.europa.eu/resource/authority/human-sex"
$terms.Hyperlinks.Add($terms.Range("B2:B7"), $url, "", "", $url) | Out-Null

# Set the friendly cell text shown for each linked cell
$terms.Range("B2").Value = "Human Sex"
$terms.Range("B3").Value = "Human Sex"
$terms.Range("B4").Value = "Human Sex"
$terms.Range("B5").Value = "Human Sex"
$terms.Range("B6").Value = "Human Sex"
$terms.Range("B7").Value = "Human Sex"

# Re-apply the Hyperlink style across the whole range (setting .Value above can reset per-cell style)
$terms.Range("B2:B7").Style = "Hyperlink"

# Restore per-sheet cursor/selection state
$wb.Worksheets.Item("fr").Range("A6").Select() | Out-Null
$wb.Worksheets.Item("pt").Range("A25").Select() | Out-Null
$wb.Worksheets.Item("si").Range("A1:E13").Select() | Out-Null
$wb.Worksheets.Item("es").Range("A1:E13").Select() | Out-Null
$wb.Worksheets.Item("en").Range("B2:D11").Select() | Out-Null
$terms.Range("B12").Select() | Out-Null

$terms.Select() | Out-Null
